# Add three new "Outliers_MAD" comparison columns (F, G, H) to the
# imputation comparison worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels in row 1, columns F:H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the header styling (bold font, thin border, centered alignment)
# used by the existing header cells (e.g. E1) by copying formats only.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Fill data rows 2-4 for the new columns with boolean FALSE values.
$ws.Range("F2:H4").Value = $false

Write-Host "Added KNN/SVM/RF Outliers_MAD columns (F:H) with header styling and FALSE values"
